$wb = $excel.ActiveWorkbook

# Mapping from short party codes (as they appear in the original workbook)
# to the new, fully-expanded party names.
$partyMap = @{
    "CDU" = "CDU - Christian Democratic Union (Christlich Demokratische Union , CDU)";
    "CSU" = "CSU - Christian Social Union  (Christlich Soziale Union, CSU)";
    "FDP" = "FDP - Free Democrats  (Freie Demokratische Partei, FDP)";
    "G -" = "G - Alliance 90-Greens (Bundnis 90-Die Grunen, G)";
    "PDS" = "PDS - Party of Democratic Socialism (Partei des Demokratischen Sozialismus, PDS)";
    "SPD" = "SPD - Social Democrats (Sozialdemokratische Partei Deutschlands, SPD)";
    "Lin" = "Linke - The Left (Die Linke, Linke)";
    "AfD" = "AfD - Alternative for Germany (Alternative für Deutschland, AfD)";
}

# Build a lookup cache so we don't recompute the same remap twice.
$cache = @{}

function Remap-Text($text) {
    if ($cache.ContainsKey($text)) {
        return $cache[$text]
    }
    # Only attempt a remap if every '+'-separated token is a known party code.
    $parts = $text.Split('+')
    $allKnown = $true
    foreach ($p in $parts) {
        if (-not $partyMap.ContainsKey($p)) {
            $allKnown = $false
            break
        }
    }
    if ($allKnown) {
        $newParts = @()
        foreach ($p in $parts) {
            $newParts += $partyMap[$p]
        }
        $joined = [string]::Join('+', $newParts)
        $cache[$text] = $joined
        return $joined
    } else {
        $cache[$text] = $text
        return $text
    }
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            $val = $cell.Value2
            if ($val -ne $null -and $val -is [string]) {
                $newVal = Remap-Text $val
                if ($newVal -ne $val) {
                    $cell.Value2 = $newVal
                }
            }
        }
    }
}
